# Update of 2025 data and RF changes
# Rows 17-59 (SpeciesSN, Code, W(kg), Numb, RF) on Sheet1 are replaced with
# the new 2025 dataset: the former RF=1 block (rows 47-59) moves up to
# rows 17-29 unchanged, and the former RF=-0.0998 block (rows 17-46) moves
# down to rows 30-59 with RF updated to 49.49275.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(17, 'Arnoglossus laterna', 'ARNOLAT', 0.08500000000000001, 7, 1),
    @(18, 'Callinectes sapidus', 'CALLSAP', 0.204, 1, 1),
    @(19, 'Chelidonichthys lucernus', 'TRIGLUC', 0.367, 4, 1),
    @(20, 'Eledone moschata', 'ELEDMOS', 0.569, 1, 1),
    @(21, 'Lithognathus mormyrus', 'LITHMOR', 1.968, 33, 1),
    @(22, 'Loligo vulgaris', 'LOLIVUL', 0.115, 1, 1),
    @(23, 'Melicertus kerathurus', 'MELIKER', 0.032, 3, 1),
    @(24, 'Ostrea edulis', 'OSTREDU', 0.03, 1, 1),
    @(25, 'Sepia officinalis', 'SEPIOFF', 1.4, 9, 1),
    @(26, 'Solea aegyptiaca', 'SOLEAEG', 0.266, 2, 1),
    @(27, 'Solea solea', 'SOLEVUL', 6.183, 68, 1),
    @(28, 'Sparus aurata', 'SPARAUR', 20.586, 211, 1),
    @(29, 'Squilla mantis', 'SQUIMAN', 2.535, 101, 1),
    @(30, 'Actiniaria nd', 'ACTINND', 0.02, 10, 49.49275),
    @(31, 'Anadara spp.', 'ANADSPP', 0.001, 2, 49.49275),
    @(32, 'Anomia ephippium', 'ANOMEPH', 0.001, 1, 49.49275),
    @(33, 'Aphrodita aculeata', 'APHRACU', 0.001, 1, 49.49275),
    @(34, 'Aporrhais pespelecani', 'APORPES', 0.038, 10, 49.49275),
    @(35, 'Ascidiacea nd', 'ASCIDND', 0.135, 23, 49.49275),
    @(36, 'Astropecten irregularis', 'ASTRIRR', 0.133, 53, 49.49275),
    @(37, 'Biological discard', 'BIOLDIS', 0.021, -1, 49.49275),
    @(38, 'Bolinus brandaris', 'MUREBRA', 0.445, 88, 49.49275),
    @(39, 'Chlamys glabra', 'CHLAGLA', 0.026, 1, 49.49275),
    @(40, 'Corbula gibba', 'CORBGIB', 0.001, 1, 49.49275),
    @(41, 'Eggs of Murex', 'EGGSMUR', 0.024, -1, 49.49275),
    @(42, 'Goneplax rhomboides', 'GONERHO', 0.018, 2, 49.49275),
    @(43, 'Gracilaria sp.', 'GRACIsp', 0.005, -1, 49.49275),
    @(44, 'Hexaplex trunculus', 'HEXATRU', 0.895, 74, 49.49275),
    @(45, 'Ilia nucleus', 'ILIANUC', 0.005, 1, 49.49275),
    @(46, 'Lima hians', 'LIMAHIA', 0.001, 1, 49.49275),
    @(47, 'Medorippe lanata', 'MEDOLAN', 0.054, 7, 49.49275),
    @(48, 'Modiolus barbatus', 'MODIBAR', 0.001, 1, 49.49275),
    @(49, 'Ocnus planci (=Cucumaria planci)', 'OCNUPLA', 0.012, 2, 49.49275),
    @(50, 'Ophiothrix sp.', 'OPHIOSP', 0.007, 11, 49.49275),
    @(51, 'Ophiura ophiura', 'OPHIOPH', 0.001, 1, 49.49275),
    @(52, 'Ostrea edulis', 'OSTREDU', 0.005, 1, 49.49275),
    @(53, 'Paguristes eremita', 'PAGUERE', 0.011, 4, 49.49275),
    @(54, 'Psammechinus microtuberculatus', 'PSAMMIC', 0.029, 9, 49.49275),
    @(55, 'Pyura dura', 'PYURADU', 0.021, 1, 49.49275),
    @(56, 'Shells NA', 'SHELLS', 0.389, -1, 49.49275),
    @(57, 'Stones NA', 'STONES', 0.866, -1, 49.49275),
    @(58, 'Ulva sp', 'ULVASPP', 0.006, -1, 49.49275),
    @(59, 'Wood NA', 'WOOD', 0.095, -1, 49.49275)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 5).Value = $row[1]   # E: SpeciesSN
    $ws.Cells.Item($r, 6).Value = $row[2]   # F: Code
    $ws.Cells.Item($r, 7).Value = $row[3]   # G: W(kg)
    $ws.Cells.Item($r, 8).Value = $row[4]   # H: Numb
    $ws.Cells.Item($r, 9).Value = $row[5]   # I: RF
}
